# Results.xlsx — rename "Sheet2" to "Volunteers characterization", move it
# right after "Volunteers", delete the now-unused "Sheet3" tab (and its
# chart), fix up the chart series that referenced the old "Sheet2" name, and
# leave the "Summary" tab active on cell V45 (matching the author's final
# on-screen state).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Repoint the charts living on "Sheet2" so their series keep working
#     once the sheet is renamed ------------------------------------------------
$charSheet = $wb.Worksheets.Item("Sheet2")
$chartObjects = $charSheet.ChartObjects()
for ($i = 1; $i -le $chartObjects.Count; $i++) {
    $chart = $chartObjects.Item($i).Chart
    $series = $chart.SeriesCollection()
    for ($j = 1; $j -le $series.Count; $j++) {
        $s = $series.Item($j)
        $s.Formula = $s.Formula.Replace("Sheet2!", "'Volunteers characterization'!")
    }
}

# --- Leave the same on-sheet selection the author had (just without the
#     frozen/scrolled topLeftCell) ------------------------------------------
$charSheet.Activate()
$charSheet.Range("H36").Select()

# --- Rename "Sheet2" to its real name and move it right after "Volunteers" --
$charSheet.Name = "Volunteers characterization"
$volunteers = $wb.Worksheets.Item("Volunteers")
$charSheet.Move($null, $volunteers)

# --- Drop the now unused "Sheet3" tab (and the chart that lived on it) -----
$wb.Worksheets.Item("Sheet3").Delete()

# --- Finish with "Summary" selected, matching the author's last action -----
$summary = $wb.Worksheets.Item("Summary")
$summary.Activate()
$summary.Range("V45").Select()
